$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply Text number format to the Price (D) cells we are about to update,
# so Excel stores the exact text (preserving formats like "38.725.61" or
# trailing zeros like "228.47" / "0.0840") instead of silently parsing them
# as floating point numbers. We restore the default "Normal" style afterward
# so no stray style index is left attached to the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '38.725.61'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '2.097.58'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '228.47'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').Value = '62.22'
$ws.Range('E7').Value = '  +1.92%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +2.01%  '
$ws.Range('D10').Value = '0.0840'
$ws.Range('E10').Value = '  -0.69%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = '15.82'
$ws.Range('E12').Value = '  +7.14%  '
$ws.Range('D13').Value = '2.408.84'
$ws.Range('D14').Value = '22.11'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = '0.803'
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '2.091.05'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').Value = '38.767.04'
$ws.Range('E18').Value = '  +1.45%  '
$ws.Range('D19').Value = '71.88'
$ws.Range('E19').Value = '  +2.15%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').Value = '227.78'
$ws.Range('E22').Value = '  +1.29%  '
$ws.Range('D24').Value = '2.37'
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('D25').Value = '2.33'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('E26').Value = '  +1.21%  '
$ws.Range('E27').Value = '  +1.99%  '
$ws.Range('E28').Value = '  +6.08%  '
$ws.Range('E29').Value = '  +4.19%  '
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('D31').Value = '2.46'
$ws.Range('E31').Value = '  +4.35%  '
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = '4.55'
$ws.Range('E33').Value = '  +2.52%  '
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('D36').Value = '6.61'
$ws.Range('E36').Value = '  +3.06%  '
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').Value = '18.28'
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.0229'
$ws.Range('E41').Value = '  +4.45%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '102.20'
$ws.Range('E42').Value = '  +2.29%  '
$ws.Range('D43').Value = '1.534.25'
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('D45').Value = '7.84'
$ws.Range('E45').Value = '  +4.62%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '1.15'
$ws.Range('E46').Value = '  +2.90%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.0911'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('E48').Value = '  -0.42%  '
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').Value = '2.293.06'
$ws.Range('E51').Value = '  +0.16%  '

# Restore default styling on the Price cells (removes the temporary text format).
$dRange.Style = "Normal"
